$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.853.61'
$ws.Range('E2').Value = '  +4.32%  '
$ws.Range('D3').Value = '2.273.47'
$ws.Range('E3').Value = '  +1.86%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '304.12'
$ws.Range('E5').Value = '  +3.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '93.21'
$ws.Range('E6').Value = '  +6.32%  '
$ws.Range('E7').Value = '  +3.54%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.487'
$ws.Range('E9').Value = '  +3.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.68'
$ws.Range('E10').Value = '  +6.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.12'
$ws.Range('E11').Value = '  +4.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0802'
$ws.Range('E12').Value = '  +2.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.115'
$ws.Range('E13').Value = '  +2.34%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.70'
$ws.Range('E14').Value = '  +3.68%  '
$ws.Range('D15').Value = '2.623.25'
$ws.Range('E15').Value = '  +2.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.26'
$ws.Range('E16').Value = '  +2.63%  '
$ws.Range('D17').Value = '2.278.80'
$ws.Range('E17').Value = '  +2.62%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.762'
$ws.Range('E18').Value = '  +3.47%  '
$ws.Range('D19').Value = '41.791.50'
$ws.Range('E19').Value = '  +4.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.31'
$ws.Range('E20').Value = '  +9.01%  '
$ws.Range('D21').Value = '0.0₃0908'
$ws.Range('E21').Value = '  +2.12%  '
$ws.Range('E22').Value = '  +2.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.38'
$ws.Range('E23').Value = '  +2.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '243.51'
$ws.Range('E24').Value = '  +2.92%  '
$ws.Range('E25').Value = '  +3.88%  '
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('E27').Value = '  +4.87%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.32'
$ws.Range('E28').Value = '  +4.24%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.63'
$ws.Range('E29').Value = '  +3.21%  '
$ws.Range('E30').Value = '  -3.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '158.18'
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '33.87'
$ws.Range('E32').Value = '  +6.51%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.20'
$ws.Range('E34').Value = '  +4.28%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0751'
$ws.Range('E35').Value = '  +4.90%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('E37').Value = '  +3.79%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '16.74'
$ws.Range('E38').Value = '  +7.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.105'
$ws.Range('E39').Value = '  +5.14%  '
$ws.Range('E40').Value = '  +3.25%  '
$ws.Range('E41').Value = '  +3.69%  '
$ws.Range('E42').Value = '  +5.70%  '
$ws.Range('D43').Value = '2.092.28'
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.53'
$ws.Range('E44').Value = '  +6.24%  '
$ws.Range('E45').Value = '  +3.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.36'
$ws.Range('E46').Value = '  +1.80%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.92'
$ws.Range('E47').Value = '  +7.15%  '
$ws.Range('E48').Value = '  +4.62%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.53'
$ws.Range('E49').Value = '  +3.33%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.11'
$ws.Range('E50').Value = '  +7.80%  '
$ws.Range('E51').Value = '  +3.47%  '
